$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the 2nd Presentation and Journals attendance columns (I, J, K) for rows 3-8
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1

$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1

$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 1

$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 1

# Update the Maximum formula in B9 (subtract 1 from the max)
$ws.Range("B9").Formula = "=MAX(B3:B8) -1"

# Update the sheet view: scroll so column G is the top-left visible column,
# and change the active selection to C12
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("C12").Select()
